# Auto-generated edit script applying cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "591.69") must be
# forced to remain text, matching the source inline-string cells, otherwise
# Excel auto-converts them to numeric cells (losing exact formatting / introducing
# floating point artifacts such as 591.69000000000005).
$textForceCells = @("D5", "D6", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D45", "D47", "D48", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.163.72"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "3.180.54"
$ws.Range("E3").Value = "  -3.84%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "591.69"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").Value = "135.33"
$ws.Range("E6").Value = "  -4.10%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.177.77"
$ws.Range("E8").Value = "  -3.77%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  -5.51%  "
$ws.Range("D11").Value = "5.23"
$ws.Range("E11").Value = "  -5.35%  "
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  -2.99%  "
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").Value = "  -3.55%  "
$ws.Range("D14").Value = "34.44"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "3.704.83"
$ws.Range("E15").Value = "  -3.80%  "
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "3.176.88"
$ws.Range("E17").Value = "  -3.86%  "
$ws.Range("D18").Value = "63.129.26"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "6.54"
$ws.Range("E19").Value = "  -4.46%  "
$ws.Range("D20").Value = "460.65"
$ws.Range("E20").Value = "  -4.05%  "
$ws.Range("D21").Value = "13.96"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "0.696"
$ws.Range("E22").Value = "  -5.49%  "
$ws.Range("D23").Value = "7.61"
$ws.Range("E23").Value = "  -4.70%  "
$ws.Range("D24").Value = "13.31"
$ws.Range("E24").Value = "  -4.28%  "
$ws.Range("D25").Value = "82.44"
$ws.Range("E25").Value = "  -3.18%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "2.67"
$ws.Range("E28").Value = "  -4.10%  "
$ws.Range("D29").Value = "7.69"
$ws.Range("E29").Value = "  -5.56%  "
$ws.Range("D30").Value = "6.74"
$ws.Range("E30").Value = "  -6.43%  "
$ws.Range("D31").Value = "2.03"
$ws.Range("E31").Value = "  -5.06%  "
$ws.Range("D32").Value = "27.29"
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("D33").Value = "0.102"
$ws.Range("E33").Value = "  -3.26%  "
$ws.Range("D34").Value = "2.37"
$ws.Range("E34").Value = "  -6.01%  "
$ws.Range("E35").Value = "  -6.05%  "
$ws.Range("D36").Value = "5.80"
$ws.Range("E36").Value = "  -3.84%  "
$ws.Range("D37").Value = "51.19"
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").Value = "0.0₃0713"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").Value = "0.0386"
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("D40").Value = "402.89"
$ws.Range("E40").Value = "  -6.33%  "
$ws.Range("D41").Value = "8.09"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("D42").Value = "2.66"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.111"
$ws.Range("E43").Value = "  -6.39%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.807.12"
$ws.Range("E44").Value = "  -10.05%  "
$ws.Range("D45").Value = "0.252"
$ws.Range("E45").Value = "  -4.25%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "2.12"
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("D48").Value = "126.27"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "25.19"
$ws.Range("E49").Value = "  -4.14%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "34.64"
$ws.Range("E50").Value = "  -5.52%  "
$ws.Range("E51").Value = "  -1.89%  "

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
